# Macroferia Regional de Talca - Haba: insert a new weekly price record.
# A new observation (2021-10-14, serial 44483) is inserted as row 16,
# pushing the existing rows 16-36 down to 17-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; Excel shifts rows 16:36 -> 17:37
# and extends the sheet dimension to A1:R37 automatically.
$ws.Rows(16).Insert()

# Populate the newly inserted row 16 with the new weekly record.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44483
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112026
$ws.Range("G16").Value = "Haba"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 8000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 8000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región de O'Higgins"
$ws.Range("P16").Value = 320
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
